$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Assigning a plain ISO-date-looking string via .Value/.Formula triggers
    # Excel's smart date-parsing heuristic and also risks allocating a brand
    # new numFmt/style record. Writing it as a quoted formula result first
    # (so it is typed as text, not parsed) and then converting the formula
    # to a static value in place via Copy + PasteSpecial(xlPasteValues)
    # keeps it a plain shared-string cell with the default (no) style -
    # exactly matching how the original rows store their date strings.
    $r = $ws.Range($addr)
    $escaped = $text.Replace('"', '""')
    $r.Formula = '="' + $escaped + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)
}

# New date strings must land in sharedStrings.xml in this exact order so the
# indices referenced by the new rows match the target workbook:
#   102 -> "2020-09-05", 103 -> "2020-09-07", 104 -> "2020-09-06"
Set-TextCell "A98" "2020-09-05"
Set-TextCell "A100" "2020-09-07"
Set-TextCell "A99" "2020-09-06"

# Row 98: 2020-09-05
$ws.Range("B98").Value = 629409
$ws.Range("C98").Value = 709182
$ws.Range("D98").Value = 86616
$ws.Range("E98").Value = 67326
$ws.Range("F98").Value = 25.06

# Row 99: 2020-09-06
$ws.Range("B99").Value = 634023
$ws.Range("C99").Value = 715395
$ws.Range("D99").Value = 82215
$ws.Range("E99").Value = 67558
$ws.Range("F99").Value = 25

# Row 100: 2020-09-07
$ws.Range("B100").Value = 637509
$ws.Range("C100").Value = 719981
$ws.Range("D100").Value = 78213
$ws.Range("E100").Value = 67781
$ws.Range("F100").Value = 24.98

# New styled, otherwise-empty row 106 (ht 18) with four cells carrying the
# new font style.
$ws.Range("I106").Font.Name = "Helvetica Neue"
$ws.Range("I106").Font.Size = 14
$ws.Range("I106").Font.Color = 3355443
$ws.Range("I106").Value = ""
$ws.Range("J106").Font.Name = "Helvetica Neue"
$ws.Range("J106").Font.Size = 14
$ws.Range("J106").Font.Color = 3355443
$ws.Range("J106").Value = ""
$ws.Range("K106").Font.Name = "Helvetica Neue"
$ws.Range("K106").Font.Size = 14
$ws.Range("K106").Font.Color = 3355443
$ws.Range("K106").Value = ""
$ws.Range("L106").Font.Name = "Helvetica Neue"
$ws.Range("L106").Font.Size = 14
$ws.Range("L106").Font.Color = 3355443
$ws.Range("L106").Value = ""
$ws.Rows("106").RowHeight = 18

$ws.Range("A87").Select()
